$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Unmerge the two trailing rows (totals row 15, footer row 16) that are
#    about to be relocated two rows further down.
# ---------------------------------------------------------------------------
$ws.Range("K15:N15").UnMerge()
$ws.Range("A16:E16").UnMerge()
$ws.Range("F16:G16").UnMerge()
$ws.Range("I16:N16").UnMerge()

# ---------------------------------------------------------------------------
# 2. Shift rows 8-16 down by two rows (bottom-up so we never clobber data that
#    still needs to move). Cut/Paste keeps each cell's existing style index
#    intact (unlike Insert/Copy, which allocate brand-new style records).
# ---------------------------------------------------------------------------
$ws.Range("A16:N16").Cut($ws.Range("A18"))
$ws.Range("A15:N15").Cut($ws.Range("A17"))
$ws.Range("A14:N14").Cut($ws.Range("A16"))
$ws.Range("A13:N13").Cut($ws.Range("A15"))
$ws.Range("A12:N12").Cut($ws.Range("A14"))
$ws.Range("A11:N11").Cut($ws.Range("A13"))
$ws.Range("A10:N10").Cut($ws.Range("A12"))
$ws.Range("A9:N9").Cut($ws.Range("A11"))
$ws.Range("A8:N8").Cut($ws.Range("A10"))

# ---------------------------------------------------------------------------
# 3. Row heights: rows keep the same alternating pattern, just shifted.
# ---------------------------------------------------------------------------
$ws.Rows(8).RowHeight = 25.5
$ws.Rows(9).RowHeight = 24.75
$ws.Rows(10).RowHeight = 25.5
$ws.Rows(11).RowHeight = 24.75
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 25.5
$ws.Rows(14).RowHeight = 24.75
$ws.Rows(15).RowHeight = 25.5
$ws.Rows(16).RowHeight = 24.75
$ws.Rows(17).RowHeight = 26.25
$ws.Rows(18).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 4. Fill in the two newly inserted product rows (8 and 9) - same layout as
#    every other data row: A=#, B=name, H=ratio, L=price, N=ratio.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = 5
$ws.Range("B8").Value2 = "HIBIOTIC 1GM 16 TAB"
$ws.Range("H8").Value2 = "2:0"
$ws.Range("L8").Value2 = 86.5
$ws.Range("N8").Value2 = "0:2"

$ws.Range("A9").Value2 = 6
$ws.Range("B9").Value2 = "KETOLAC 10MG 20 TAB"
$ws.Range("H9").Value2 = "0:1"
$ws.Range("L9").Value2 = 19
$ws.Range("N9").Value2 = "0:2"

# ---------------------------------------------------------------------------
# 5. Re-number the "#" column for every row that shifted (items now at rows
#    10-16 are items 7-13).
# ---------------------------------------------------------------------------
$ws.Range("A10").Value2 = 7
$ws.Range("A11").Value2 = 8
$ws.Range("A12").Value2 = 9
$ws.Range("A13").Value2 = 10
$ws.Range("A14").Value2 = 11
$ws.Range("A15").Value2 = 12
$ws.Range("A16").Value2 = 13

# ---------------------------------------------------------------------------
# 6. Re-merge: new data rows 15 & 16 need the usual B:G / H:K / L:M merges;
#    the relocated totals/footer rows need their merges restored too.
# ---------------------------------------------------------------------------
$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("B16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()

$ws.Range("K17:N17").Merge()
$ws.Range("A18:E18").Merge()
$ws.Range("F18:G18").Merge()
$ws.Range("I18:N18").Merge()

# ---------------------------------------------------------------------------
# 7. Update the grand-total cell for the new balance.
# ---------------------------------------------------------------------------
$ws.Range("K17").Value2 = 413.5
